# Redeem points 76442781 500.0
# Appends a redemption record row (phone, points, timestamp) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: phone number, kept as literal text (it's all digits, so without
# forcing a text format it would be stored as a number). Briefly apply a
# text NumberFormat so the value is accepted as a string, then strip the
# formatting back off so the cell keeps the workbook's default style.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "76442781"
$ws.Range("A2").ClearFormats()

# B2: points redeemed, a plain number.
$ws.Range("B2").Value = 500

# C2: ISO-ish timestamp string, also kept as literal text.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025-08-20T08:55:01"
$ws.Range("C2").ClearFormats()
